$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrialData")

# Rename the three date-related headers (sow_date/harvest_date/em_date -> date_sow/date_harvest/date_emerg)
$ws.Range("N1").Value = "date_sow"
$ws.Range("O1").Value = "date_harvest"
$ws.Range("P1").Value = "date_emerg"

# Correct the emergence date for trial row 2 (one day earlier)
$ws.Range("P2").Value = 42118

# Re-enter the id formula across A2:A8 so Excel collapses it into one shared formula
$ws.Range("A2:A8").Formula = "=C2*10000+E2+G2*100-200000"

# Update the view: scroll so column E is at the left edge and select P3
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P3").Select() | Out-Null
